$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("reaction19")

# Clear the entire row first (removes C1:Q1 entirely) then set the new values
$ws.Range("A1:Q1").ClearContents()

$ws.Range("A1").Value = 38
$ws.Range("B1").Value = 39
